# Updates cryptos list with latest prices/volume figures
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Cells whose new value is a plain number-looking string must be
# forced to Text format first so Excel stores them as text (matching
# the original inline-string cell type) instead of auto-converting
# them into numeric values.
$textCells = @("D5","D6","D12","D15","D21","D22","D28","D35","D36","D38","D39","D41","D44","D45","D46","D47","D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "69.470.88"
$ws.Range("E2").Value = "  -2.63%  "
$ws.Range("D3").Value = "3.693.01"
$ws.Range("E3").Value = "  -3.10%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "688.99"
$ws.Range("E5").Value = "  -1.68%  "
$ws.Range("D6").Value = "161.99"
$ws.Range("D7").Value = "3.691.86"
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  -5.51%  "
$ws.Range("E10").Value = "  -8.08%  "
$ws.Range("E11").Value = "  -1.80%  "
$ws.Range("D12").Value = "0.439"
$ws.Range("E12").Value = "  -8.51%  "
$ws.Range("E13").Value = "  -5.77%  "
$ws.Range("D14").Value = "4.314.60"
$ws.Range("E14").Value = "  -3.14%  "
$ws.Range("D15").Value = "33.07"
$ws.Range("E15").Value = "  -8.16%  "
$ws.Range("D16").Value = "3.689.73"
$ws.Range("E16").Value = "  -2.90%  "
$ws.Range("D17").Value = "69.465.05"
$ws.Range("E17").Value = "  -2.68%  "
$ws.Range("E18").Value = "  -1.25%  "
$ws.Range("E19").Value = "  -8.27%  "
$ws.Range("E20").Value = "  -9.62%  "
$ws.Range("D21").Value = "476.74"
$ws.Range("E21").Value = "  -7.40%  "
$ws.Range("D22").Value = "9.96"
$ws.Range("E22").Value = "  -4.93%  "
$ws.Range("E23").Value = "  -7.91%  "
$ws.Range("E24").Value = "  -4.88%  "
$ws.Range("D25").Value = "3.836.96"
$ws.Range("E25").Value = "  -3.08%  "
$ws.Range("E26").Value = "  -9.19%  "
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("D28").Value = "11.23"
$ws.Range("E28").Value = "  -10.04%  "
$ws.Range("E29").Value = "  -10.47%  "
$ws.Range("E30").Value = "  -11.22%  "
$ws.Range("E31").Value = "  -9.98%  "
$ws.Range("E32").Value = "  -7.73%  "
$ws.Range("E33").Value = "  -8.09%  "
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").Value = "0.166"
$ws.Range("E35").Value = "  -3.77%  "
$ws.Range("B36").Value = "EthereumClassic"
$ws.Range("C36").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D36").Value = "26.88"
$ws.Range("E36").Value = "  -7.75%  "
$ws.Range("D37").Value = "3.656.28"
$ws.Range("E37").Value = "  -3.12%  "
$ws.Range("D38").Value = "8.35"
$ws.Range("E38").Value = "  -9.16%  "
$ws.Range("D39").Value = "6.29"
$ws.Range("E39").Value = "  -2.18%  "
$ws.Range("E40").Value = "  -4.71%  "
$ws.Range("D41").Value = "0.0918"
$ws.Range("E41").Value = "  -8.99%  "
$ws.Range("E43").Value = "  -0.05%  "
$ws.Range("D44").Value = "0.953"
$ws.Range("E44").Value = "  -6.28%  "
$ws.Range("D45").Value = "163.53"
$ws.Range("E45").Value = "  -4.94%  "
$ws.Range("D46").Value = "48.33"
$ws.Range("E46").Value = "  -2.91%  "
$ws.Range("D47").Value = "30.30"
$ws.Range("E48").Value = "  -15.30%  "
$ws.Range("E49").Value = "  -8.42%  "
$ws.Range("D50").Value = "1.32"
$ws.Range("E50").Value = "  -3.84%  "
$ws.Range("E51").Value = "  -2.41%  "
